$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19 ("fuselagem" row), shifting rows 19-29 down to 20-30
$ws.Rows("19:19").Insert()

# New column D is used for long annotation text; widen it like the author did
$ws.Columns("D").ColumnWidth = 16.140625

$ws.Range("G4").Value = "x"
$ws.Range("H4").Value = "m"
$ws.Range("F5").Value = 2
$ws.Range("G5").Formula = "=2.175-0.5"
$ws.Range("H5").Value = 200
$ws.Range("I5").Formula = "=G5*H5"
$ws.Range("M5").Value = "pos asa/ fuselagem"
$ws.Range("O5").Value = 0.2
$ws.Range("R5").Value = "centro fuselagem"
$ws.Range("F6").Value = 2
$ws.Range("G6").Formula = "=G5+0.81"
$ws.Range("H6").Value = 200
$ws.Range("I6").Formula = "=G6*H6"
$ws.Range("R6").Value = "x"
$ws.Range("S6").Value = "A"
$ws.Range("F7").Value = 2
$ws.Range("G7").Formula = "=G6+0.81"
$ws.Range("H7").Value = 200
$ws.Range("I7").Formula = "=G7*H7"
$ws.Range("L7").Formula = "=G7-G6"
$ws.Range("Q7").Value = "cone 1"
$ws.Range("R7").Value = 0.67
$ws.Range("S7").Formula = "=1*2.25*0.5"
$ws.Range("T7").Formula = "=R7*S7"
$ws.Range("F8").Value = 2
$ws.Range("G8").Formula = "=G7+0.81"
$ws.Range("H8").Value = 200
$ws.Range("I8").Formula = "=G8*H8"
$ws.Range("Q8").Value = "centro"
$ws.Range("R8").Formula = "=1+6.5/2"
$ws.Range("S8").Formula = "=2.25*6.5"
$ws.Range("T8").Formula = "=R8*S8"
$ws.Range("F9").Value = 2
$ws.Range("G9").Formula = "=6.025"
$ws.Range("H9").Value = 200
$ws.Range("I9").Formula = "=G9*H9"
$ws.Range("K9").Value = "45cm entre `"parede`" da roda e banco"
$ws.Range("Q9").Value = "cone 2"
$ws.Range("R9").Formula = "=6.5+1.5/3"
$ws.Range("S9").Formula = "=(2.25*1.5)/2"
$ws.Range("T9").Formula = "=R9*S9"
$ws.Range("F10").Value = "1+rt"
$ws.Range("G10").Formula = "=7.025"
$ws.Range("H10").Formula = "=80+100"
$ws.Range("I10").Formula = "=G10*H10"
$ws.Range("K10").Value = "temos 80 cm para a parte das rodas"
$ws.Range("S10").Formula = "=SUM(S7:S9)"
$ws.Range("F11").Value = "bagagens"
$ws.Range("G11").Formula = "=2.1"
$ws.Range("H11").Value = 200
$ws.Range("I11").Formula = "=G11*H11"
$ws.Range("F12").Value = "mot f"
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 140
$ws.Range("I12").Formula = "=G12*H12"
$ws.Range("S12").Value = "xcfus"
$ws.Range("T12").Formula = "=SUM(T7:T9)/S10"
$ws.Range("F13").Value = "mot tras"
$ws.Range("G13").Value = 6.5
$ws.Range("H13").Value = 140
$ws.Range("I13").Formula = "=G13*H13"
$ws.Range("F14").Value = "Tbgerad"
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 460
$ws.Range("I14").Formula = "=G14*H14"
$ws.Range("F15").Value = "bat+av"
$ws.Range("G15").Value = 0.5
$ws.Range("H15").Value = 510
$ws.Range("I15").Formula = "=G15*H15"
$ws.Range("F16").Value = "rf"
$ws.Range("G16").Formula = "=1-0.28/2"
$ws.Range("H16").Value = 50
$ws.Range("I16").Formula = "=G16*H16"
$ws.Range("O16").Value = "para ft com x=2.5"
$ws.Range("F17").Value = "asa frente"
$ws.Range("G17").Formula = "=G12"
$ws.Range("H17").Value = 200
$ws.Range("I17").Formula = "=G17*H17"
$ws.Range("O17").Value = "pos asa fuselagem = 0.7"
$ws.Range("F18").Value = "asa trás"
$ws.Range("G18").Formula = "=G13"
$ws.Range("H18").Value = 200
$ws.Range("I18").Formula = "=G18*H18"
$ws.Range("F19").Value = "fuselagem"
$ws.Range("G19").Formula = "=T12"
$ws.Range("H19").Value = 800
$ws.Range("I19").Formula = "=G19*H19"
$ws.Range("F20").Value = "ft"
$ws.Range("G20").Formula = "=2.25"
$ws.Range("H20").Value = 2200
$ws.Range("I20").Formula = "=G20*H20"
$ws.Range("K20").Value = "valor muito exagerado???"
$ws.Range("H21").Formula = "=SUM(H5:H20)"
$ws.Range("I21").Formula = "=SUM(I5:I20)"
$ws.Range("K21").Value = "xcm"
$ws.Range("K22").Formula = "=I21/H21"
$ws.Range("E24").Value = "c"
$ws.Range("F24").Value = 2
$ws.Range("I24").Value = "hn"
$ws.Range("J24").Formula = "=(F27 + (F28/F24)*(1-F29))/(1+(1-F29))"
$ws.Range("D25").Value = "posi. 1ª asa"
$ws.Range("E25").Value = "wp"
$ws.Range("F25").Formula = "=G17-1"
$ws.Range("I25").Value = "xpn"
$ws.Range("J25").Formula = "=F25+J24*F24"
$ws.Range("K25").Value = "m"
$ws.Range("L25").Formula = "=J25-K22"
$ws.Range("E26").Value = "aw"
$ws.Range("F26").Formula = "=0.106"
$ws.Range("G26").Value = "[/grau]"
$ws.Range("E27").Value = "h_nw"
$ws.Range("F27").Value = 0.25
$ws.Range("I27").Value = "SM"
$ws.Range("J27").Formula = "=(J25-K22)/(2*F24)"
$ws.Range("D28").Value = "considerado pelo prof como distancia entre bordo da primeira asa e ponto neutro da segunda"
$ws.Range("E28").Value = "lht"
$ws.Range("F28").Formula = "=G18-G17+F24*F27"
$ws.Range("E29").Value = "e_alfa"
$ws.Range("F29").Value = 0.4
$ws.Range("I29").Value = "Cl_alfa"
$ws.Range("J29").Formula = "=F26+F26*(1-F29)"
$ws.Range("K29").Value = "[/grau]"
$ws.Range("I30").Value = "cm_alfa"
$ws.Range("J30").Formula = "=-J27*J29"
$ws.Range("K30").Value = "[/grau]"

# Match the final selection left by the author
$ws.Range("G20").Select() | Out-Null
